$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Rushing")
$ws2 = $wb.Worksheets.Item("Receiving")

foreach ($ws in @($ws1, $ws2)) {
    # Insert a new column at I, shifting existing player columns right.
    $ws.Columns("I:I").Insert()

    # Copy the header formatting (bold font + border + centered) from the
    # neighboring header cell so the new header matches the others.
    $ws.Range("H1").Copy()
    $ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
    $excel.CutCopyMode = $false

    # New player column: header name + placeholder "n" data value.
    $ws.Range("I1").Value = "De.Jackson"
    $ws.Range("I2").Value = "n"
}
